$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 44/45 coin name & link swap (VeChain <-> Maker)
$ws.Range('B44').Value = 'Maker'
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# Price column (D) updates - force text to preserve formatting (e.g. trailing zeros, multi-dot numbers)
Set-TextValue $ws.Range('D2') '60.887.28'
Set-TextValue $ws.Range('D3') '2.917.88'
Set-TextValue $ws.Range('D5') '590.50'
Set-TextValue $ws.Range('D6') '146.50'
Set-TextValue $ws.Range('D13') '33.59'
Set-TextValue $ws.Range('D15') '3.401.10'
Set-TextValue $ws.Range('D16') '60.809.77'
Set-TextValue $ws.Range('D18') '2.919.45'
Set-TextValue $ws.Range('D19') '430.29'
Set-TextValue $ws.Range('D20') '13.37'
Set-TextValue $ws.Range('D21') '0.678'
Set-TextValue $ws.Range('D23') '81.44'
Set-TextValue $ws.Range('D24') '10.91'
Set-TextValue $ws.Range('D26') '11.84'
Set-TextValue $ws.Range('D30') '7.01'
Set-TextValue $ws.Range('D31') '26.61'
Set-TextValue $ws.Range('D32') '0.109'
Set-TextValue $ws.Range('D33') '1.00'
Set-TextValue $ws.Range('D36') '5.62'
Set-TextValue $ws.Range('D40') '8.56'
Set-TextValue $ws.Range('D41') '0.283'
Set-TextValue $ws.Range('D42') '40.22'
Set-TextValue $ws.Range('D43') '379.73'
Set-TextValue $ws.Range('D44') '2.694.41'
Set-TextValue $ws.Range('D45') '0.0343'
Set-TextValue $ws.Range('D46') '133.31'
Set-TextValue $ws.Range('D48') '23.74'
Set-TextValue $ws.Range('D50') '2.00'

# Volume(1h) column (E) updates
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  +5.63%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -2.84%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('E39').Value = '  -4.01%  '
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('E51').Value = '  -0.27%  '

